$wb = $excel.ActiveWorkbook
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
$ws = $wb.Worksheets.Item("Ver-Iniciação1")
Write-Host "Found sheet: " $ws.Name
Write-Host "C2 value before: " $ws.Range("C2").Value
$ws.Range("C2").Value = 44474
Write-Host "C2 value after: " $ws.Range("C2").Value
Write-Host "C2 NumberFormat before: " $ws.Range("C2").NumberFormat
$ws.Range("C2").NumberFormat = "m/d/yyyy"
Write-Host "C2 NumberFormat after: " $ws.Range("C2").NumberFormat
Write-Host "F2 formula: " $ws.Range("F2").Formula
Write-Host "F2 value: " $ws.Range("F2").Value
